$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the "Implement edge detector(Canny?)" run into three runs with
#    proofErr (gramStart/gramEnd) markers bracketing "detector(".
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Implement edge detector(Canny?)") {
        $target = $i
        break
    }
}

$p = $d.Paragraphs.Item($target)
$r = $p.Range
# Range excluding the trailing paragraph mark so InsertXML only replaces the
# runs, leaving the paragraph's own pPr (style/numbering) untouched.
$contentRange = $d.Range($r.Start, $r.End - 1)

$splitXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Implement edge </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>detector(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Canny?)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$contentRange.InsertXML($splitXml)

# ---------------------------------------------------------------------------
# 2) Insert two new list-items after that paragraph:
#      "Blur image to remove detail"
#      "Convert image to grayscale"  (yellow highlight; this is also where
#                                      the "_GoBack" bookmark now lives,
#                                      having moved on from its old home at
#                                      the end of "Manually enter correct
#                                      classification")
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item($target)
$p.Range.InsertParagraphAfter() | Out-Null
$pBlur = $d.Paragraphs.Item($target + 1)
$pBlur.Range.InsertParagraphAfter() | Out-Null

$pBlur = $d.Paragraphs.Item($target + 1)
$blurXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Blur image to remove detail</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rngBlur = $d.Range($pBlur.Range.Start, $pBlur.Range.End)
$rngBlur.InsertXML($blurXml)

$pGray = $d.Paragraphs.Item($target + 2)
$grayXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Convert image to grayscale</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rngGray = $d.Range($pGray.Range.Start, $pGray.Range.End)
$rngGray.InsertXML($grayXml)
